$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values (row 1 headers, row 2 data) ---
$ws.Range("H1").Value = "SocialSecurityBenefits"
$ws.Range("I1").Value = "maritalStatus"
$ws.Range("J1").Value = "social_security_override"
$ws.Range("H2").Value = "Yes"
$ws.Range("I2").Value = "Married"
$ws.Range("J2").Value = 14576

# --- Column widths for the newly added columns ---
$ws.Columns.Item(8).ColumnWidth = 31.44140625
$ws.Columns.Item(9).ColumnWidth = 22.6640625
$ws.Columns.Item(10).ColumnWidth = 35.5546875

# --- Make H1 / J1 start from the same font/format as E1 (big font, no explicit color,
#     vertical-centered) before we touch alignment, so they end up sharing E1's style ---
$ws.Range("E1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

# --- Re-apply the values since PasteSpecial(formats) shouldn't clobber them, but make sure ---
$ws.Range("H1").Value = "SocialSecurityBenefits"
$ws.Range("J1").Value = "social_security_override"

# --- Header row formatting: center the plain header cells (A:D, F:G) ---
$ws.Range("A1:G1").HorizontalAlignment = -4108

# --- H1 / J1 match E1's special combined alignment (horizontal + vertical centered) ---
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4108

# --- I1 gets its own new font (size 14, explicit black color) + centered alignment ---
$ws.Range("I1").Font.Size = 14
$ws.Range("I1").Font.Color = 0
$ws.Range("I1").HorizontalAlignment = -4108

# --- F2 / G2 drop their old distinct style and match the rest of row 2 ---
$ws.Range("F2:G2").Font.Size = 12

# --- View: scroll right a bit and move the selection ---
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K6").Select()
